$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 228.89473
$ws.Range("I9").Value = 44.842106
$ws.Range("J9").Value = 412.94736
$ws.Range("K9").Value = 44.842106
$ws.Range("L9").Value = 412.94736
$ws.Range("M9").Value = 124.157894
$ws.Range("N9").Value = -750.94736

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3222.6155
$ws.Range("I64").Value = 3025.5715
$ws.Range("J64").Value = 3265.7188
$ws.Range("K64").Value = 3025.5715
$ws.Range("L64").Value = 3265.7188
$ws.Range("M64").Value = -2777.5715
$ws.Range("N64").Value = -3761.7188

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3222.6155
$ws.Range("I67").Value = 3025.5715
$ws.Range("J67").Value = 3265.7188
$ws.Range("K67").Value = 3025.5715
$ws.Range("L67").Value = 3265.7188
$ws.Range("M67").Value = -2167.5715
$ws.Range("N67").Value = -4981.718800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 27299.41
$ws.Range("I129").Value = 747.5
$ws.Range("J129").Value = 45770.305
$ws.Range("K129").Value = 2242.5
$ws.Range("L129").Value = 137310.915
$ws.Range("M129").Value = 2757.5
$ws.Range("N129").Value = -147310.915

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1808.8889
$ws.Range("I137").Value = 1293.3334
$ws.Range("J137").Value = 2066.6667
$ws.Range("K137").Value = 3880.0002
$ws.Range("L137").Value = 6200.000100000001
$ws.Range("M137").Value = -1330.0002
$ws.Range("N137").Value = -11300.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2419.6445
$ws.Range("I138").Value = 778.375
$ws.Range("J138").Value = 3732.66
$ws.Range("K138").Value = 2335.125
$ws.Range("L138").Value = 11197.98
$ws.Range("M138").Value = 2804.875
$ws.Range("N138").Value = -21477.98

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2490.5334
$ws.Range("I122").Value = 2490.5334
$ws.Range("K122").Value = 7471.600199999999
$ws.Range("M122").Value = -5021.600199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1197604.4
$ws.Range("I105").Value = 1624665.9
$ws.Range("J105").Value = 1832.2
$ws.Range("K105").Value = 1624665.9
$ws.Range("L105").Value = 1832.2
$ws.Range("M105").Value = -1622918.9
$ws.Range("N105").Value = -5326.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2748.9333
$ws.Range("I134").Value = 2710.3076
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 8130.9228
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -5595.9228
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2475.25
$ws.Range("I31").Value = 1340.591
$ws.Range("J31").Value = 6635.6665
$ws.Range("K31").Value = 1340.591
$ws.Range("L31").Value = 6635.6665
$ws.Range("M31").Value = -1045.591
$ws.Range("N31").Value = -7225.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2475.25
$ws.Range("I34").Value = 1340.591
$ws.Range("J34").Value = 6635.6665
$ws.Range("K34").Value = 1340.591
$ws.Range("L34").Value = 6635.6665
$ws.Range("M34").Value = -1138.591
$ws.Range("N34").Value = -7039.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1254.2858
$ws.Range("I105").Value = 884.44446
$ws.Range("J105").Value = 1920
$ws.Range("K105").Value = 884.44446
$ws.Range("L105").Value = 1920
$ws.Range("M105").Value = 862.55554
$ws.Range("N105").Value = -5414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1000.9286
$ws.Range("I122").Value = 911
$ws.Range("J122").Value = 1162.8
$ws.Range("K122").Value = 2733
$ws.Range("L122").Value = 3488.4
$ws.Range("M122").Value = -283
$ws.Range("N122").Value = -8388.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6000
$ws.Range("I56").Value = 6000
$ws.Range("K56").Value = 6000
$ws.Range("M56").Value = -5470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 740.41174
$ws.Range("J122").Value = 948.3
$ws.Range("L122").Value = 8534.699999999999
$ws.Range("N122").Value = -13434.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 85.5
$ws.Range("I2").Value = 87
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 87
$ws.Range("L2").Value = 75
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = -301

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2958.375
$ws.Range("I102").Value = 2334.5293
$ws.Range("K102").Value = 2334.5293
$ws.Range("M102").Value = -712.5293000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1514.2
$ws.Range("I122").Value = 1642.75
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4928.25
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2478.25
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -95070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2111.889
$ws.Range("I7").Value = 1943.8572
$ws.Range("J7").Value = 2700
$ws.Range("K7").Value = 1943.8572
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = -1831.8572
$ws.Range("N7").Value = -2924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2012.5
$ws.Range("I32").Value = 1371.4286
$ws.Range("J32").Value = 6500
$ws.Range("K32").Value = 1371.4286
$ws.Range("L32").Value = 6500
$ws.Range("M32").Value = -1054.4286
$ws.Range("N32").Value = -7134

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2939.7856
$ws.Range("I40").Value = 2673.2273
$ws.Range("K40").Value = 2673.2273
$ws.Range("M40").Value = -2537.2273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1640.2
$ws.Range("I46").Value = 1433.6666
$ws.Range("J46").Value = 1950
$ws.Range("K46").Value = 1433.6666
$ws.Range("L46").Value = 1950
$ws.Range("M46").Value = -1245.6666
$ws.Range("N46").Value = -2326

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3938.862
$ws.Range("I122").Value = 3216
$ws.Range("J122").Value = 8456.75
$ws.Range("K122").Value = 9648
$ws.Range("L122").Value = 25370.25
$ws.Range("M122").Value = -7198
$ws.Range("N122").Value = -30270.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2111.889
$ws.Range("I126").Value = 1943.8572
$ws.Range("J126").Value = 2700
$ws.Range("K126").Value = 5831.571599999999
$ws.Range("L126").Value = 8100
$ws.Range("M126").Value = -3361.571599999999
$ws.Range("N126").Value = -13040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2512.7334
$ws.Range("I136").Value = 2062.818
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 6188.454000000001
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -3638.454000000001
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 6790.8335
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 7749
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 7749
$ws.Range("M45").Value = -1509
$ws.Range("N45").Value = -8731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1886.4706
$ws.Range("I122").Value = 1390.8334
$ws.Range("J122").Value = 3076
$ws.Range("K122").Value = 4172.5002
$ws.Range("L122").Value = 9228
$ws.Range("M122").Value = -1722.5002
$ws.Range("N122").Value = -14128

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1344.4445
$ws.Range("I126").Value = 1166.6666
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 3499.9998
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -1029.9998
$ws.Range("N126").Value = -10040
